$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24+ down by one.
$ws.Rows("24:24").Insert()

# The insert copies formatting from the row above (which used the
# word-wrap style in column B); clear it so the new row uses the default
# (unstyled) formatting, matching the original row 24 it is replacing.
$ws.Rows("24:24").ClearFormats()

# Populate the newly inserted row 24 with the new key/value pair.
$ws.Range("A24").Value = "health_warning"
$ws.Range("B24").Value = "Watch out! Once the heart bar is empty, you will have to start over!"

# Move the selection to the newly inserted row, matching the saved view state.
$ws.Range("B24").Select()
